$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.957.58"
$ws.Range("D3").Value = "2.677.03"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.26"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.57"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "2.676.61"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  -4.72%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").Value = "3.170.55"
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").Value = "71.886.89"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.22"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "2.675.90"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.19"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.69"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.98"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.77"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "0.0₃0972"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "501.91"
$ws.Range("E32").Value = "  -7.78%  "
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.65"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.63"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.06"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.00"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.333"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.51"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.89"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.607"
$ws.Range("E51").Value = "  +0.03%  "
